# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with freshly scraped values (GitHub Actions refresh).
# Note: some Price values look like plain numbers (e.g. "94.70"); they are
# prefixed with a leading apostrophe so Excel stores/keeps them as literal
# text (matching the original inlineStr cell content) instead of silently
# converting them to floating point numbers and losing formatting/precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.408.45"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.288.62"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "'94.70"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("D7").Value = "'0.505"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").Value = "'34.29"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "'18.87"
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "2.647.47"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "2.285.78"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "'0.775"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "42.352.44"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "'12.13"
$ws.Range("E19").Value = "  -5.69%  "
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'5.97"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "'67.40"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "'234.99"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +5.31%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "'24.11"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "'2.36"
$ws.Range("E28").Value = "  +14.97%  "
$ws.Range("D29").Value = "'163.95"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").Value = "'8.99"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").Value = "'31.52"
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'17.35"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").Value = "'0.0691"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'2.32"
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("D37").Value = "'4.34"
$ws.Range("E37").Value = "  -8.01%  "
$ws.Range("D38").Value = "'0.0992"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D41").Value = "'2.67"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "'19.81"
$ws.Range("E42").Value = "  +9.25%  "
$ws.Range("D43").Value = "1.944.59"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("D45").Value = "'0.0276"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("E46").Value = "  +2.71%  "
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").Value = "2.517.32"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "'2.82"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").Value = "'52.59"
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("E51").Value = "  +0.41%  "
